# Update cryptos list with refreshed price / volume(1h) figures,
# and reorder rows 48-49 (THORChain now ranks above FLOKI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as text (many entries use "." as a
# thousands separator, e.g. "71.745.41", which is not a valid number).
# Force the whole column to Text format first so plain-looking numeric
# strings (e.g. "519.21", "1.00") are not silently coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.271.97"
$ws.Range("E2").Value = "  +4.23%  "
$ws.Range("D3").Value = "4.042.75"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "519.21"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "148.21"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "0.738"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  +2.58%  "
$ws.Range("D11").Value = "0.0000335"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "47.74"
$ws.Range("E12").Value = "  +13.39%  "
$ws.Range("D13").Value = "10.82"
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").Value = "4.691.14"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "4.049.37"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "21.22"
$ws.Range("E16").Value = "  +7.43%  "
$ws.Range("D17").Value = "14.23"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "72.220.76"
$ws.Range("E20").Value = "  +4.34%  "
$ws.Range("D21").Value = "438.72"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "97.41"
$ws.Range("E22").Value = "  +11.08%  "
$ws.Range("D23").Value = "3.52"
$ws.Range("E23").Value = "  +5.79%  "
$ws.Range("D24").Value = "14.52"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "11.94"
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "11.23"
$ws.Range("E27").Value = "  +5.98%  "
$ws.Range("D28").Value = "36.99"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").Value = "3.07"
$ws.Range("E29").Value = "  +9.37%  "
$ws.Range("D30").Value = "704.73"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "13.46"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").Value = "7.18"
$ws.Range("E32").Value = "  +21.68%  "
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("D34").Value = "68.16"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").Value = "0.0₃0900"
$ws.Range("E35").Value = "  +8.13%  "
$ws.Range("D36").Value = "3.72"
$ws.Range("E36").Value = "  +26.65%  "
$ws.Range("D37").Value = "0.436"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "40.56"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("E39").Value = "  +4.00%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.0486"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("D44").Value = "2.75"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "3.53"
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("D46").Value = "0.145"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("D50").Value = "3.33"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "0.0₆0339"
$ws.Range("E51").Value = "  +2.14%  "

# Rows 48 and 49 swap coin order: THORChain moves to rank 46 (row 48),
# FLOKI moves to rank 47 (row 49). Rank numbers in column A are unchanged.
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "9.04"
$ws.Range("E48").Value = "  +6.97%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000272"
$ws.Range("E49").Value = "  +21.99%  "
